# Atualização no script de automação
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 3 through 8 (old comprovante entries no longer needed)
$ws.Range("A3:E8").EntireRow.Delete()

# Update remaining data row with the latest comprovante info
$ws.Range("A2").Value = "Você"
$ws.Range("B2").Value = "Desconhecido"
$ws.Range("C2").Value = "R$ 67,45"
$ws.Range("D2").Value = "SHIBATA COMERCIO E ATACADO DE PROD..."
$ws.Range("E2").Value = "Motoboy"

# Adjust column widths (stored XML width = ColumnWidth + 5/6 in this runtime)
$ws.Columns.Item(1).ColumnWidth = 5.166666666666667
$ws.Columns.Item(2).ColumnWidth = 13.166666666666666
$ws.Columns.Item(3).ColumnWidth = 9.166666666666666
